$d = $word.ActiveDocument
$script:bmCounter = 0

function Split-At($pos) {
    # Forces a run boundary at $pos without altering any text. Uses a
    # throwaway bookmark: adding it (to a collapsed range) breaks the run
    # the position falls inside into two runs with identical (inherited)
    # formatting; deleting the bookmark afterwards leaves the split intact
    # but removes the bookmark markup itself.
    $r = $d.Range($pos, $pos)
    $bmName = "zzsplit" + $script:bmCounter
    $script:bmCounter = $script:bmCounter + 1
    $d.Bookmarks.Add($bmName, $r)
    $d.Bookmarks($bmName).Delete()
}

function Insert-SplitText($pos, $txt) {
    # Inserts $txt at $pos as its own standalone run, split off from
    # whatever runs precede/follow it.
    $r = $d.Range($pos, $pos)
    $r.InsertBefore($txt)
    Split-At $pos
    $endPos = $pos + $txt.Length
    Split-At $endPos
}

function Replace-SplitText($start, $end, $newTxt) {
    # Replaces the text in [start, end) with $newTxt, leaving $newTxt as
    # its own standalone run, split off from its neighbours.
    $r = $d.Range($start, $end)
    $r.Text = $newTxt
    Split-At $start
    $endPos = $start + $newTxt.Length
    Split-At $endPos
}

# ---------------------------------------------------------------------------
# 1) "... ударом команды будет передан ..." -> "... ударом команды, будет
#    передан ..." (insert a comma right after "команды").
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "ударом команды будет"
$idx = $text.IndexOf($needle)
$prefixLen = "ударом команды".Length
$pos = $idx + $prefixLen
Insert-SplitText $pos ","

# ---------------------------------------------------------------------------
# 2) Drop the _GoBack bookmark that used to sit between "... состояние " and
#    "наказывается ...".
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 3) "тайм аут" -> "тайм-аут" (the space between becomes a hyphen) inside
#    "Каждая команда может 2 раза за игру тайм аут длиной ...".
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "Каждая команда может 2 раза за игру тайм аут"
$idx = $text.IndexOf($needle)
$prefixLen = "Каждая команда может 2 раза за игру тайм".Length
$start = $idx + $prefixLen
$end = $start + 1
Replace-SplitText $start $end "-"

# ---------------------------------------------------------------------------
# 4) "Разрашеатся" -> "Разрешеатся" (typo fix: the 5th letter "а" -> "е").
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "Разрашеатся"
$idx = $text.IndexOf($needle)
$prefixLen = "Разр".Length
$start = $idx + $prefixLen
$end = $start + 1
Replace-SplitText $start $end "е"

# ---------------------------------------------------------------------------
# 5) "... для подачи на линию защиты ..." -> "... для ввода в игру на линию
#    защиты ...".
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "Наказание за задержку – передача мяча противнику для подачи на линию"
$idx = $text.IndexOf($needle)
$prefixLen = "Наказание за задержку – передача мяча противнику для ".Length
$start = $idx + $prefixLen
$oldLen = "подачи".Length
$end = $start + $oldLen
Replace-SplitText $start $end "ввода в игру"

# ---------------------------------------------------------------------------
# 6) "... противником в первых двух случаях. Далее ..." -> "... противником
#    первые два раза. Далее ...".
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "Наказание – подача мяча противником в первых двух случаях"
$idx = $text.IndexOf($needle)
$prefixLen = "Наказание – подача мяча противником ".Length
$start = $idx + $prefixLen
$oldLen = "в первых двух случаях".Length
$end = $start + $oldLen
Replace-SplitText $start $end "первые два раза"

# ---------------------------------------------------------------------------
# 7) "... владении) мешающих ... владения збрасывается." ->
#    "... владении), мешающих ... владения сбрасывается."
#    (insert a comma after the closing paren, and fix "збрасывается" ->
#    "сбрасывается").
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "оставаясь в вашем владении) мешающих"
$idx = $text.IndexOf($needle)
$prefixLen = "оставаясь в вашем владении)".Length
$pos = $idx + $prefixLen
Insert-SplitText $pos ","

$text = $d.Content.Text
$needle = "время владения збрасывается"
$idx = $text.IndexOf($needle)
$prefixLen = "время владения ".Length
$start = $idx + $prefixLen
$end = $start + 1
Replace-SplitText $start $end "с"

# ---------------------------------------------------------------------------
# 8) "... противников. Наказанием  – подача ..." -> "... противников.
#    Наказание – подача ..." (drop trailing "м"), and re-add the _GoBack
#    bookmark right after the new "Наказание" run.
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "Запрещены без согласия судьи или противников. Наказанием "
$idx = $text.IndexOf($needle)
$prefixLen = "Запрещены без согласия судьи или противников. Наказание".Length
$start = $idx + $prefixLen
$end = $start + 1
Replace-SplitText $start $end ""

$text = $d.Content.Text
$needle2 = "Запрещены без согласия судьи или противников. Наказание"
$idx2 = $text.IndexOf($needle2)
$pos2 = $idx2 + $needle2.Length
$bmRange = $d.Range($pos2, $pos2)
$d.Bookmarks.Add("_GoBack", $bmRange)
